$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe in the assigned text makes Excel store the value
# as literal text (exactly like typing an apostrophe before a number in the
# UI) without changing the cells NumberFormat. This keeps numeric-looking
# strings such as "1.00" / "5.80" / "3.00" intact with their trailing zeros,
# matching the original inlineStr cells in the workbook.

$ws.Range("D2").Value = '''68.229.98'
$ws.Range("E2").Value = '''  +1.36%  '
$ws.Range("D3").Value = '''3.903.91'
$ws.Range("E3").Value = '''  +0.77%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '''  +0.11%  '
$ws.Range("D5").Value = '''480.37'
$ws.Range("E5").Value = '''  +2.14%  '
$ws.Range("D6").Value = '''144.94'
$ws.Range("E6").Value = '''  -0.30%  '
$ws.Range("D7").Value = '''0.622'
$ws.Range("E7").Value = '''  -2.02%  '
$ws.Range("E8").Value = '''  -0.06%  '
$ws.Range("D9").Value = '''0.723'
$ws.Range("E9").Value = '''  -3.12%  '
$ws.Range("E10").Value = '''  +6.64%  '
$ws.Range("E11").Value = '''  +12.17%  '
$ws.Range("D12").Value = '''42.58'
$ws.Range("E12").Value = '''  -2.10%  '
$ws.Range("D13").Value = '''10.56'
$ws.Range("E13").Value = '''  +0.89%  '
$ws.Range("D14").Value = '''4.532.48'
$ws.Range("E14").Value = '''  +0.70%  '
$ws.Range("D15").Value = '''14.58'
$ws.Range("E15").Value = '''  -1.64%  '
$ws.Range("D16").Value = '''3.893.37'
$ws.Range("E16").Value = '''  +0.67%  '
$ws.Range("E17").Value = '''  -0.43%  '
$ws.Range("D18").Value = '''19.68'
$ws.Range("E18").Value = '''  -2.10%  '
$ws.Range("E19").Value = '''  -3.39%  '
$ws.Range("D20").Value = '''68.325.04'
$ws.Range("E20").Value = '''  +1.17%  '
$ws.Range("D21").Value = '''436.36'
$ws.Range("E21").Value = '''  +0.02%  '
$ws.Range("D22").Value = '''3.36'
$ws.Range("E22").Value = '''  +1.47%  '
$ws.Range("E23").Value = '''  -2.41%  '
$ws.Range("D24").Value = '''87.69'
$ws.Range("E24").Value = '''  -1.82%  '
$ws.Range("E25").Value = '''  +17.81%  '
$ws.Range("E26").Value = '''  -1.18%  '
$ws.Range("D27").Value = '''38.07'
$ws.Range("E27").Value = '''  -0.33%  '
$ws.Range("D28").Value = '''10.38'
$ws.Range("E28").Value = '''  +3.02%  '
$ws.Range("D29").Value = '''5.80'
$ws.Range("E29").Value = '''  +5.51%  '
$ws.Range("D30").Value = '''702.60'
$ws.Range("E30").Value = '''  -3.70%  '
$ws.Range("B31").Value = '''Hedera'
$ws.Range("C31").Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = '''0.131'
$ws.Range("E31").Value = '''  -2.38%  '
$ws.Range("B32").Value = '''Cosmos'
$ws.Range("C32").Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").Value = '''13.33'
$ws.Range("E32").Value = '''  -3.91%  '
$ws.Range("D33").Value = '''2.86'
$ws.Range("E33").Value = '''  +2.82%  '
$ws.Range("D34").Value = '''0.0₃0920'
$ws.Range("E34").Value = '''  +34.64%  '
$ws.Range("D35").Value = '''41.50'
$ws.Range("E35").Value = '''  -6.46%  '
$ws.Range("D36").Value = '''59.31'
$ws.Range("E36").Value = '''  +2.03%  '
$ws.Range("D37").Value = '''5.73'
$ws.Range("E37").Value = '''  +3.25%  '
$ws.Range("E38").Value = '''  -7.79%  '
$ws.Range("E39").Value = '''  -0.16%  '
$ws.Range("E40").Value = '''  -2.39%  '
$ws.Range("E41").Value = '''  +10.85%  '
$ws.Range("D42").Value = '''2.76'
$ws.Range("E42").Value = '''  +7.84%  '
$ws.Range("D43").Value = '''3.00'
$ws.Range("E43").Value = '''  +1.98%  '
$ws.Range("E44").Value = '''  -2.17%  '
$ws.Range("E45").Value = '''  -0.41%  '
$ws.Range("E46").Value = '''  -0.06%  '
$ws.Range("D47").Value = '''3.43'
$ws.Range("E47").Value = '''  -1.38%  '
$ws.Range("E48").Value = '''  -0.89%  '
$ws.Range("D49").Value = '''145.95'
$ws.Range("E49").Value = '''  +1.01%  '
$ws.Range("E50").Value = '''  -4.71%  '
$ws.Range("E51").Value = '''  -2.20%  '
